# Updates cryptos list values (prices / volume %) per Dec 15 2023 commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.176.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.248.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.633"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.56%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.04"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0957"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.102"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.580.09"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.83"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.246.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.038.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0985"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.11"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.39"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.63"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.84%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +13.04%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.60"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0829"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.91"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.39"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0303"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.90"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.19"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "113.81"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.204"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.28"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.34%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.446"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +17.20%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.17"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.49%  "
